$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "last updated" timestamp footer (cell A1) ---
$ws.Range("A1").Value = "Datos actualizados a 8 de Agosto de 2020 a las 19:30"

# --- Country name re-ranking (rows whose rank order changed) ---
$ws.Range("A35").Value = "Israel"
$ws.Range("A36").Value = "Oman"
$ws.Range("A69").Value = "Etiopia"
$ws.Range("A70").Value = "Costa Rica"
$ws.Range("A71").Value = "Austria"
$ws.Range("A149").Value = "Siria"
$ws.Range("A150").Value = "Gambia"
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Updated daily statistics (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes) ---
$ws.Range("B4").Value = 5115470
$ws.Range("C4").Value = 19946
$ws.Range("D4").Value = 2619317
$ws.Range("E4").Value = 2331643
$ws.Range("G4").Value = 416
$ws.Range("H4").Value = 164510
$ws.Range("B5").Value = 2988796
$ws.Range("C5").Value = 21732
$ws.Range("E5").Value = 820162
$ws.Range("G5").Value = 538
$ws.Range("H5").Value = 100240
$ws.Range("B6").Value = 2150858
$ws.Range("C6").Value = 63994
$ws.Range("D6").Value = 1476994
$ws.Range("E6").Value = 630418
$ws.Range("G6").Value = 868
$ws.Range("H6").Value = 43446
$ws.Range("B11").Value = 371023
$ws.Range("C11").Value = 2198
$ws.Range("D11").Value = 344133
$ws.Range("E11").Value = 16879
$ws.Range("G11").Value = 53
$ws.Range("H11").Value = 10011
$ws.Range("B31").Value = 93572
$ws.Range("C31").Value = 1603
$ws.Range("D31").Value = 71605
$ws.Range("E31").Value = 16051
$ws.Range("G31").Value = 19
$ws.Range("H31").Value = 5916
$ws.Range("B35").Value = 82279
$ws.Range("C35").Value = 1288
$ws.Range("D35").Value = 57068
$ws.Range("E35").Value = 24619
$ws.Range("G35").Value = 11
$ws.Range("H35").Value = 592
$ws.Range("B36").Value = 81357
$ws.Range("C36").Value = 290
$ws.Range("D36").Value = 73481
$ws.Range("E36").Value = 7367
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 509
$ws.Range("B59").Value = 34639
$ws.Range("C59").Value = 484
$ws.Range("D59").Value = 24083
$ws.Range("E59").Value = 9263
$ws.Range("G59").Value = 11
$ws.Range("H59").Value = 1293
$ws.Range("B61").Value = 32007
$ws.Range("C61").Value = 1345
$ws.Range("D61").Value = 22190
$ws.Range("E61").Value = 9337
$ws.Range("G61").Value = 19
$ws.Range("H61").Value = 480
$ws.Range("B65").Value = 26644
$ws.Range("C65").Value = 174
$ws.Range("E65").Value = 1508
$ws.Range("B66").Value = 25837
$ws.Range("C66").Value = 699
$ws.Range("D66").Value = 11899
$ws.Range("E66").Value = 13520
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 418
$ws.Range("B69").Value = 22253
$ws.Range("C69").Value = 801
$ws.Range("D69").Value = 9707
$ws.Range("E69").Value = 12156
$ws.Range("G69").Value = 10
$ws.Range("H69").Value = 390
$ws.Range("B70").Value = 22081
$ws.Range("C70").Value = 0
$ws.Range("D70").Value = 7266
$ws.Range("E70").Value = 14597
$ws.Range("G70").Value = 0
$ws.Range("H70").Value = 218
$ws.Range("B71").Value = 21919
$ws.Range("C71").Value = 82
$ws.Range("D71").Value = 19812
$ws.Range("E71").Value = 1386
$ws.Range("G71").Value = 1
$ws.Range("H71").Value = 721
$ws.Range("E79").Value = 6126
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 96
$ws.Range("B84").Value = 11754
$ws.Range("C84").Value = 200
$ws.Range("D84").Value = 7622
$ws.Range("E84").Value = 3609
$ws.Range("G84").Value = 4
$ws.Range("H84").Value = 523
$ws.Range("B125").Value = 2463
$ws.Range("C125").Value = 13
$ws.Range("E125").Value = 1241
$ws.Range("B149").Value = 1125
$ws.Range("C149").Value = 65
$ws.Range("D149").Value = 331
$ws.Range("E149").Value = 744
$ws.Range("G149").Value = 2
$ws.Range("H149").Value = 50
$ws.Range("B150").Value = 1090
$ws.Range("D150").Value = 146
$ws.Range("E150").Value = 925
$ws.Range("H150").Value = 19
